# Generate Report for Handoff
#
# Refresh the "Latest Handoff Datetime" (column D) for every row whose
# handoff just completed (i.e. every data row except the ones that already
# have a full target/handback pair filled in, and except the sentinel
# ".localization-config" row). Each locale sheet gets its own single new
# timestamp, reused across all of the rows it touches - mirroring how the
# CI report-generator stamps one "now" per locale run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhRows = @(7, 10, 12, 13, 14, 15, 16)
$deRows = @(7, 10, 12, 13, 14, 15, 16)

foreach ($r in $zhRows) {
    $zhcn.Range("D$r").Value = "2016-03-09 10:29:49"
}

foreach ($r in $deRows) {
    $dede.Range("D$r").Value = "2016-03-09 10:29:54"
}
